$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "cache policy" answer cell (B9) already carries the target style
# (s=3). Copy its format onto the other answer cells in column B so the
# new "Y" values end up styled the same way Excel would normalize them.
$ws.Range("B9").Copy()
$ws.Range("B2:B8").PasteSpecial(-4122)
$ws.Range("B10:B15").PasteSpecial(-4122)

$ws.Range("B2:B8").Value = "Y"
$ws.Range("B10:B15").Value = "Y"

# Rows 16-18 previously had no B cell; writing a value there picks up
# column B's default style (s=4) automatically.
$ws.Range("B16:B18").Value = "Y"

$ws.Range("B19").Select()
